$wb = $excel.ActiveWorkbook

# --- Create both new sheets first (tab order: Day 22 then Day 23) ---
$wsNew22 = $wb.Worksheets.Add()
$wsNew22.Name = "Day 22"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew22.Move($null, $lastSheet)

$wsNew23 = $wb.Worksheets.Add()
$wsNew23.Name = "Day 23"
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew23.Move($null, $lastSheet2)

# References taken before .Move() go stale and further writes against them
# are silently dropped, so re-fetch both sheets by name now that tab order
# is settled.
$ws22 = $wb.Worksheets.Item("Day 22")
$ws23 = $wb.Worksheets.Item("Day 23")

# --- Populate "Day 23" (teacher_id, subject_id, dept_id) FIRST so its
#     header strings land earlier in the shared-string table, matching
#     the workbook's original authoring order. ---
$ws23.Range("A1").Value = "teacher_id"
$ws23.Range("B1").Value = "subject_id"
$ws23.Range("C1").Value = "dept_id"

$day23Data = @(
    @(1, 2, 3),
    @(1, 2, 4),
    @(1, 3, 3),
    @(2, 1, 1),
    @(2, 2, 1),
    @(2, 3, 1),
    @(2, 4, 1)
)

$r = 2
foreach ($row in $day23Data) {
    $ws23.Cells.Item($r, 1).Value = $row[0]
    $ws23.Cells.Item($r, 2).Value = $row[1]
    $ws23.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# --- Populate "Day 22" (player_id, device_id, event_date, games_played) ---
$ws22.Range("A1").Value = "player_id"
$ws22.Range("B1").Value = "device_id"
$ws22.Range("C1").Value = "event_date"
$ws22.Range("D1").Value = "games_played"

$day22Data = @(
    @(1, 2, 42430, 5),
    @(1, 2, 42492, 6),
    @(2, 3, 42911, 1),
    @(3, 1, 42431, 0),
    @(3, 4, 43284, 5)
)

$r = 2
foreach ($row in $day22Data) {
    $ws22.Cells.Item($r, 1).Value = $row[0]
    $ws22.Cells.Item($r, 2).Value = $row[1]
    $ws22.Cells.Item($r, 3).Value = $row[2]
    $ws22.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Apply the workbook's existing built-in date format (m/d/yyyy, style index
# 2) to the event_date column by copy/paste-special of formats from an
# existing date cell elsewhere in the workbook. Setting .NumberFormat
# directly always mints a brand-new custom numFmt entry instead of reusing
# the built-in one already present in styles.xml, which would diverge from
# the original file's style table.
$dateFormatSource = $wb.Worksheets.Item("Day 21").Range("B2")
$dateFormatSource.Copy()
$ws22.Range("C2:C6").PasteSpecial(-4122)

# --- Make "Day 22" the active/selected sheet ---
$wb.Worksheets.Item("Day 22").Activate()
